$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text / shared-string cells (mirrors existing A:G header/label cells into I:O) ---
$ws.Range("I2").Value = "Week 2"
$ws.Range("J2").Value = "A01"
$ws.Range("K4").Value = "Documentation"
$ws.Range("J5").Value = "Måndag"
$ws.Range("K5").Value = "Tisdag"
$ws.Range("L5").Value = "Onsdag"
$ws.Range("M5").Value = "Torsdag"
$ws.Range("N5").Value = "Fredag"
$ws.Range("I6").Value = "Calle"
$ws.Range("I7").Value = "Kim"
$ws.Range("I8").Value = "Nils"
$ws.Range("I9").Value = "Rasmus"
$ws.Range("K11").Value = "Strategy meeting"
$ws.Range("J12").Value = "Måndag"
$ws.Range("K12").Value = "Tisdag"
$ws.Range("L12").Value = "Onsdag"
$ws.Range("M12").Value = "Torsdag"
$ws.Range("N12").Value = "Fredag"
$ws.Range("I13").Value = "Calle"
$ws.Range("I14").Value = "Kim"
$ws.Range("I15").Value = "Nils"
$ws.Range("I16").Value = "Rasmus"
$ws.Range("K18").Value = "Conceptual view"
$ws.Range("J19").Value = "Måndag"
$ws.Range("K19").Value = "Tisdag"
$ws.Range("L19").Value = "Onsdag"
$ws.Range("M19").Value = "Torsdag"
$ws.Range("N19").Value = "Fredag"
$ws.Range("I20").Value = "Calle"
$ws.Range("I21").Value = "Kim"
$ws.Range("I22").Value = "Nils"
$ws.Range("I23").Value = "Rasmus"
$ws.Range("K25").Value = "Assignment planning"
$ws.Range("J26").Value = "Måndag"
$ws.Range("K26").Value = "Tisdag"
$ws.Range("L26").Value = "Onsdag"
$ws.Range("M26").Value = "Torsdag"
$ws.Range("N26").Value = "Fredag"
$ws.Range("I27").Value = "Calle"
$ws.Range("I28").Value = "Kim"
$ws.Range("I29").Value = "Nils"
$ws.Range("I30").Value = "Rasmus"
$ws.Range("C34").Value = "Factors and Issues"
$ws.Range("K34").Value = "Factors and Issues"
$ws.Range("B35").Value = "Måndag"
$ws.Range("C35").Value = "Tisdag"
$ws.Range("D35").Value = "Onsdag"
$ws.Range("E35").Value = "Torsdag"
$ws.Range("F35").Value = "Fredag"
$ws.Range("J35").Value = "Måndag"
$ws.Range("K35").Value = "Tisdag"
$ws.Range("L35").Value = "Onsdag"
$ws.Range("M35").Value = "Torsdag"
$ws.Range("N35").Value = "Fredag"
$ws.Range("A36").Value = "Calle"
$ws.Range("I36").Value = "Calle"
$ws.Range("A37").Value = "Kim"
$ws.Range("I37").Value = "Kim"
$ws.Range("A38").Value = "Nils"
$ws.Range("I38").Value = "Nils"
$ws.Range("A39").Value = "Rasmus"
$ws.Range("I39").Value = "Rasmus"

# --- Formula cells (mirrors existing SUM formulas into the new J:O block) ---
$ws.Range("O6").Formula = "=SUM(J6:N6)"
$ws.Range("O7").Formula = "=SUM(J7:N7)"
$ws.Range("O8").Formula = "=SUM(J8:N8)"
$ws.Range("O9").Formula = "=SUM(J9:N9)"
$ws.Range("O10").Formula = "=SUM(O6:O9)"
$ws.Range("O13").Formula = "=SUM(J13:N13)"
$ws.Range("O14").Formula = "=SUM(J14:N14)"
$ws.Range("O15").Formula = "=SUM(J15:N15)"
$ws.Range("O16").Formula = "=SUM(J16:N16)"
$ws.Range("O17").Formula = "=SUM(O13:O16)"
$ws.Range("O20").Formula = "=SUM(J20:N20)"
$ws.Range("O21").Formula = "=SUM(J21:N21)"
$ws.Range("O22").Formula = "=SUM(J22:N22)"
$ws.Range("O23").Formula = "=SUM(J23:N23)"
$ws.Range("O24").Formula = "=SUM(O20:O23)"
$ws.Range("O27").Formula = "=SUM(J27:N27)"
$ws.Range("O28").Formula = "=SUM(J28:N28)"
$ws.Range("O29").Formula = "=SUM(J29:N29)"
$ws.Range("O30").Formula = "=SUM(J30:N30)"
$ws.Range("O31").Formula = "=SUM(O27:O30)"
$ws.Range("G36").Formula = "=SUM(B36:F36)"
$ws.Range("O36").Formula = "=SUM(J36:N36)"
$ws.Range("G37").Formula = "=SUM(B37:F37)"
$ws.Range("O37").Formula = "=SUM(J37:N37)"
$ws.Range("G38").Formula = "=SUM(B38:F38)"
$ws.Range("O38").Formula = "=SUM(J38:N38)"
$ws.Range("G39").Formula = "=SUM(B39:F39)"
$ws.Range("O39").Formula = "=SUM(J39:N39)"
$ws.Range("G40").Formula = "=SUM(G36:G39)"
$ws.Range("O40").Formula = "=SUM(O36:O39)"

# --- Numeric input cells (K36:K39 = 4, highlighted like the other entered hour values) ---
$ws.Range("K36").Value = 4
$ws.Range("K37").Value = 4
$ws.Range("K38").Value = 4
$ws.Range("K39").Value = 4
$numAddrs = @("K36","K37","K38","K39")
foreach ($a in $numAddrs) { $ws.Range($a).Interior.Color = 15986394 }

# --- Blank cells carrying an explicit "no fill" style (matches s=3 in the target) ---
$blankAddrs = @("J7","J27","J28","J29","J30","J31","B36","C36","J36","B37","C37","J37","B38","C38","J38","B39","C39","J39","B40","J40")
foreach ($a in $blankAddrs) { $ws.Range($a).Interior.ColorIndex = -4142 }

# --- Totals-row fill (matches s=1 used by the existing G10/G17/G24/G31 totals) ---
$totalAddrs = @("O10","O17","O24","O31","G40","O40")
foreach ($a in $totalAddrs) { $ws.Range($a).Interior.Color = 5880731 }

# --- Sheet view: scroll position + final selection ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C41").Select()
